$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 35-37 (Pseudotime_1, Pseudotime_2, Pseudotime_3), shifting rows below up.
$ws.Range("A35:D37").EntireRow.Delete()
